$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.009723801316837
$ws.Range("D2").Value = 1.012369409256422
$ws.Range("E2").Value = 1.012102512257772
$ws.Range("F2").Value = 1.013057754868212
$ws.Range("I2").Value = 1.022892173712927
$ws.Range("J2").Value = 1.014982178193998
$ws.Range("K2").Value = 1.01523377473692
$ws.Range("L2").Value = 1.014967681395946
$ws.Range("M2").Value = 1.015920049790344
$ws.Range("N2").Value = 1.00906612447811

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.0123290077605
$ws.Range("D3").Value = 1.014885506124189
$ws.Range("E3").Value = 1.014349969949083
$ws.Range("F3").Value = 1.01657539013031
$ws.Range("I3").Value = 1.023246752324957
$ws.Range("J3").Value = 1.017211050148072
$ws.Range("K3").Value = 1.017549598742093
$ws.Range("L3").Value = 1.017015553002236
$ws.Range("M3").Value = 1.019234790840363
$ws.Range("N3").Value = 1.009804260830041

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.013990547240399
$ws.Range("D4").Value = 1.016489922865572
$ws.Range("E4").Value = 1.015783641761536
$ws.Range("F4").Value = 1.018789415519919
$ws.Range("I4").Value = 1.023461817648929
$ws.Range("J4").Value = 1.018629430864997
$ws.Range("K4").Value = 1.019024248906082
$ws.Range("L4").Value = 1.018319830112994
$ws.Range("M4").Value = 1.021317697615991
$ws.Range("N4").Value = 1.010273859763054

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.014683415119066
$ws.Range("D5").Value = 1.017158900337124
$ws.Range("E5").Value = 1.016381557804524
$ws.Range("F5").Value = 1.019705639075996
$ws.Range("I5").Value = 1.023548832358737
$ws.Range("J5").Value = 1.019220147637942
$ws.Range("K5").Value = 1.01963862512403
$ws.Range("L5").Value = 1.018863284695991
$ws.Range("M5").Value = 1.022178827605462
$ws.Range("N5").Value = 1.010469404293567

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.014799424001475
$ws.Range("D6").Value = 1.017270905055041
$ws.Range("E6").Value = 1.016481672689962
$ws.Range("F6").Value = 1.019858632571337
$ws.Range("I6").Value = 1.023563244356673
$ws.Range("J6").Value = 1.019319008826121
$ws.Range("K6").Value = 1.019741459127039
$ws.Range("L6").Value = 1.018954251391702
$ws.Range("M6").Value = 1.022322572007372
$ws.Range("N6").Value = 1.010502128436787

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.013999827354155
$ws.Range("D7").Value = 1.016498883273337
$ws.Range("E7").Value = 1.015791649838458
$ws.Range("F7").Value = 1.018801714916446
$ws.Range("I7").Value = 1.023462993645269
$ws.Range("J7").Value = 1.018637345760307
$ws.Range("K7").Value = 1.019032479922084
$ws.Range("L7").Value = 1.018327110730025
$ws.Range("M7").Value = 1.021329260751614
$ws.Range("N7").Value = 1.010276479945826

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.010609347934248
$ws.Range("D8").Value = 1.013224726849681
$ws.Range("E8").Value = 1.012866392416505
$ws.Range("F8").Value = 1.014259616781333
$ws.Range("I8").Value = 1.02301500687937
$ws.Range("J8").Value = 1.015740459935344
$ws.Range("K8").Value = 1.016021442638626
$ws.Range("L8").Value = 1.015664159477811
$ws.Range("M8").Value = 1.017053300973298
$ws.Range("N8").Value = 1.009317271314984

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.004442324129462
$ws.Range("D9").Value = 1.007267006087641
$ws.Range("E9").Value = 1.007547956565901
$ws.Range("F9").Value = 1.005764727953794
$ws.Range("I9").Value = 1.022113488764371
$ws.Range("J9").Value = 1.010446490180462
$ws.Range("K9").Value = 1.010526254529018
$ws.Range("L9").Value = 1.010806230397623
$ws.Range("M9").Value = 1.009029197930934
$ws.Range("N9").Value = 1.007563355763422

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.000190790630187
$ws.Range("D10").Value = 1.003158226520186
$ws.Range("E10").Value = 1.003883126305623
$ws.Range("F10").Value = 0.9997481347662661
$ws.Range("I10").Value = 1.021433968023802
$ws.Range("J10").Value = 1.006780010644065
$ws.Range("K10").Value = 1.006725420328503
$ws.Range("L10").Value = 1.007447541029766
$ws.Range("M10").Value = 1.003328462018642
$ws.Range("N10").Value = 1.006347980336754

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 0.9983139923082324
$ws.Range("D11").Value = 1.001344077935216
$ws.Range("E11").Value = 1.002265739850661
$ws.Range("F11").Value = 0.9970534203241391
$ws.Range("I11").Value = 1.021120325440467
$ws.Range("J11").Value = 1.005157452936121
$ws.Range("K11").Value = 1.005044597867223
$ws.Range("L11").Value = 1.00596257883049
$ws.Range("M11").Value = 1.00077117612243
$ws.Range("N11").Value = 1.00580997892648

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 0.9976112416549122
$ws.Range("D12").Value = 1.000664730304424
$ws.Range("E12").Value = 1.001660187812595
$ws.Range("F12").Value = 0.9960385151339523
$ws.Range("I12").Value = 1.021000835876873
$ws.Range("J12").Value = 1.004549292025317
$ws.Range("K12").Value = 1.004414778155151
$ws.Range("L12").Value = 1.005406200460119
$ws.Range("M12").Value = 0.9998074338895888
$ws.Range("N12").Value = 1.005608304957111

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 0.9977622425398768
$ws.Range("D13").Value = 1.000810705075609
$ws.Range("E13").Value = 1.001790300592335
$ws.Range("F13").Value = 0.9962568566521686
$ws.Range("I13").Value = 1.021026603280786
$ws.Range("J13").Value = 1.004679995863504
$ws.Range("K13").Value = 1.004550128691058
$ws.Range("L13").Value = 1.005525765883156
$ws.Range("M13").Value = 1.000014795257649
$ws.Range("N13").Value = 1.005651649031043

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 0.9982560188774017
$ws.Range("D14").Value = 1.001288036271832
$ws.Range("E14").Value = 1.00221578350893
$ws.Range("F14").Value = 0.9969698160516492
$ws.Range("I14").Value = 1.021110509809414
$ws.Range("J14").Value = 1.005107295091406
$ws.Range("K14").Value = 1.004992650059767
$ws.Range("L14").Value = 1.005916687461878
$ws.Range("M14").Value = 1.000691798532269
$ws.Range("N14").Value = 1.005793346392718

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 0.9985594981624113
$ws.Range("D15").Value = 1.001581400868412
$ws.Range("E15").Value = 1.002477297565917
$ws.Range("F15").Value = 0.9974072266498927
$ws.Range("I15").Value = 1.021161809064245
$ws.Range("J15").Value = 1.005369836491277
$ws.Range("K15").Value = 1.005264568008201
$ws.Range("L15").Value = 1.006156905386268
$ws.Range("M15").Value = 1.001107071245156
$ws.Range("N15").Value = 1.005880405219693

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.000314572203269
$ws.Range("D16").Value = 1.003277868402726
$ws.Range("E16").Value = 1.003989807583235
$ws.Range("F16").Value = 0.999925044806819
$ws.Range("I16").Value = 1.021454368471693
$ws.Range("J16").Value = 1.006886939675252
$ws.Range("K16").Value = 1.006836214051649
$ws.Range("L16").Value = 1.007545431365877
$ws.Range("M16").Value = 1.003496266692218
$ws.Range("N16").Value = 1.006383432325287

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.001405721331596
$ws.Range("D17").Value = 1.004332483604167
$ws.Range("E17").Value = 1.004930263503624
$ws.Range("F17").Value = 1.001480088691472
$ws.Range("I17").Value = 1.021632637890683
$ws.Range("J17").Value = 1.007829070916578
$ws.Range("K17").Value = 1.007812532590323
$ws.Range("L17").Value = 1.008408082874746
$ws.Range("M17").Value = 1.004970817198845
$ws.Range("N17").Value = 1.006695775789255

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.002038727365172
$ws.Range("D18").Value = 1.004944260258981
$ws.Range("E18").Value = 1.005475887925492
$ws.Range("F18").Value = 1.002378522313599
$ws.Range("I18").Value = 1.021734752321862
$ws.Range("J18").Value = 1.008375243796897
$ws.Range("K18").Value = 1.008378638062154
$ws.Range("L18").Value = 1.008908311687686
$ws.Range("M18").Value = 1.005822362926067
$ws.Range("N18").Value = 1.006876833297565

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.002253989060768
$ws.Range("D19").Value = 1.00515229670407
$ws.Range("E19").Value = 1.00566144109552
$ws.Range("F19").Value = 1.002683422020291
$ws.Range("I19").Value = 1.021769256141534
$ws.Range("J19").Value = 1.008560911793445
$ws.Range("K19").Value = 1.008571101253682
$ws.Range("L19").Value = 1.009078383643104
$ws.Range("M19").Value = 1.006111285187604
$ws.Range("N19").Value = 1.006938380187218

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.001289009126888
$ws.Range("D20").Value = 1.004219682798374
$ws.Range("E20").Value = 1.004829665716956
$ws.Range("F20").Value = 1.001314140109574
$ws.Range("I20").Value = 1.021613704903941
$ws.Range("J20").Value = 1.007728337857735
$ws.Range("K20").Value = 1.007708132438668
$ws.Range("L20").Value = 1.008315834130888
$ws.Range("M20").Value = 1.004813498209291
$ws.Range("N20").Value = 1.006662381391307

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 0.998110771290399
$ws.Range("D21").Value = 1.001147627668387
$ws.Range("E21").Value = 1.002090623107367
$ws.Range("F21").Value = 0.9967602572629009
$ws.Range("I21").Value = 1.02108588454948
$ws.Range("J21").Value = 1.004981618979655
$ws.Range("K21").Value = 1.004862491895571
$ws.Range("L21").Value = 1.005801704866674
$ws.Range("M21").Value = 1.000492824588544
$ws.Range("N21").Value = 1.005751671353364

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 0.9960798183746542
$ws.Range("D22").Value = 0.9991842038725113
$ws.Range("E22").Value = 1.000340697474055
$ws.Range("F22").Value = 0.9938159473963067
$ws.Range("I22").Value = 1.020736690443559
$ws.Range("J22").Value = 1.003222876024386
$ws.Range("K22").Value = 1.003041454659543
$ws.Range("L22").Value = 1.004193109875636
$ws.Range("M22").Value = 0.9976958220205044
$ws.Range("N22").Value = 1.005168407839271

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 0.9971596448138133
$ws.Range("D23").Value = 1.000228156824833
$ws.Range("E23").Value = 1.001271070320574
$ws.Range("F23").Value = 0.9953846540837086
$ws.Range("I23").Value = 1.020923473761725
$ws.Range("J23").Value = 1.004158307918429
$ws.Range("K23").Value = 1.004009920431604
$ws.Range("L23").Value = 1.00504856632192
$ws.Range("M23").Value = 0.9991863680966246
$ws.Range("N23").Value = 1.005478643367967

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.001341756971019
$ws.Range("D24").Value = 1.004270663002171
$ws.Range("E24").Value = 1.004875130559628
$ws.Range("F24").Value = 1.001389151660494
$ws.Range("I24").Value = 1.021622265670234
$ws.Range("J24").Value = 1.007773865147567
$ws.Range("K24").Value = 1.00775531675527
$ws.Range("L24").Value = 1.008357526442902
$ws.Range("M24").Value = 1.004884610226809
$ws.Range("N24").Value = 1.006677474360367

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.006060528904564
$ws.Range("D25").Value = 1.008830557187946
$ws.Range("E25").Value = 1.008943207516789
$ws.Range("F25").Value = 1.008021340258952
$ws.Range("I25").Value = 1.022360144424155
$ws.Range("J25").Value = 1.011838519317112
$ws.Range("K25").Value = 1.011970327927236
$ws.Range("L25").Value = 1.008024655735144
$ws.Range("M25").Value = 1.011163807811549
$ws.Range("N25").Value = 1.008024655735144
